$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row - add Commission column
$ws.Range("G1").Value = "Commission"

# Update existing rows 2 and 3 with new data, replacing old rows
$ws.Range("A2").Value = "26-10-2022"
$ws.Range("B2").Value = "Buy"
$ws.Range("C2").Value = "EREGL.IS"
$ws.Range("D2").Value = 30.28
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = 302.8
$ws.Range("G2").Value = 0.63

$ws.Range("A3").Value = "27-10-2022"
$ws.Range("B3").Value = "Buy"
$ws.Range("C3").Value = "EREGL.IS"
$ws.Range("D3").Value = 30.86
$ws.Range("E3").Value = 15
$ws.Range("F3").Value = 462.9
$ws.Range("G3").Value = 0.97

$ws.Range("A4").Value = "08-11-2022"
$ws.Range("B4").Value = "Buy"
$ws.Range("C4").Value = "EREGL.IS"
$ws.Range("D4").Value = 36
$ws.Range("E4").Value = 10
$ws.Range("F4").Value = 360
$ws.Range("G4").Value = 0.75

$ws.Range("A5").Value = "-"
$ws.Range("B5").Value = "Total"
$ws.Range("C5").Value = "EREGL.IS"
$ws.Range("D5").Value = 32.163
$ws.Range("E5").Value = 35
$ws.Range("F5").Value = 1125.7
$ws.Range("G5").Value = 2.35
